$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("python_libraries")
Write-Host $ws.Name
Write-Host $ws.Range("A1").Value
Write-Host $ws.Range("A2").Value
Write-Host $ws.Range("B2").Value
